$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1179.8  # H32: 1034.4 -> 1179.8
$ws.Cells.Item(32, 9).Value = 0  # I32: 990 -> 0
$ws.Cells.Item(32, 10).Value = 1179.8  # J32: 1078.8 -> 1179.8
$ws.Cells.Item(32, 11).Value = 0  # K32: 990 -> 0
$ws.Cells.Item(32, 12).Value = 1179.8  # L32: 1078.8 -> 1179.8
$ws.Cells.Item(32, 13).ClearContents()  # M32: -664 -> (removed)
$ws.Cells.Item(32, 14).Value = -1831.8  # N32: -1730.8 -> -1831.8

$ws.Cells.Item(41, 8).Value = 373.69232  # H41: 378.125 -> 373.69232
$ws.Cells.Item(41, 9).Value = 476.85715  # I41: 453 -> 476.85715
$ws.Cells.Item(41, 11).Value = 476.85715  # K41: 453 -> 476.85715
$ws.Cells.Item(41, 13).Value = -36.85714999999999  # M41: -13 -> -36.85714999999999

$ws.Cells.Item(62, 8).Value = 6833  # H62: 5731.25 -> 6833
$ws.Cells.Item(62, 9).Value = 7168.3335  # I62: 4435 -> 7168.3335
$ws.Cells.Item(62, 10).Value = 6665.3335  # J62: 6163.3335 -> 6665.3335
$ws.Cells.Item(62, 11).Value = 7168.3335  # K62: 4435 -> 7168.3335
$ws.Cells.Item(62, 12).Value = 6665.3335  # L62: 6163.3335 -> 6665.3335
$ws.Cells.Item(62, 13).Value = -6544.3335  # M62: -3811 -> -6544.3335
$ws.Cells.Item(62, 14).Value = -7913.3335  # N62: -7411.3335 -> -7913.3335

$ws.Cells.Item(65, 8).Value = 6833  # H65: 5731.25 -> 6833
$ws.Cells.Item(65, 9).Value = 7168.3335  # I65: 4435 -> 7168.3335
$ws.Cells.Item(65, 10).Value = 6665.3335  # J65: 6163.3335 -> 6665.3335
$ws.Cells.Item(65, 11).Value = 35841.6675  # K65: 22175 -> 35841.6675
$ws.Cells.Item(65, 12).Value = 33326.6675  # L65: 30816.6675 -> 33326.6675
$ws.Cells.Item(65, 13).Value = -32721.6675  # M65: -19055 -> -32721.6675
$ws.Cells.Item(65, 14).Value = -39566.6675  # N65: -37056.6675 -> -39566.6675

$ws.Cells.Item(76, 8).Value = 3083.2778  # H76: 3145.4546 -> 3083.2778
$ws.Cells.Item(76, 9).Value = 3046.6  # I76: 3114.2856 -> 3046.6
$ws.Cells.Item(76, 10).Value = 3266.6667  # J76: 3200 -> 3266.6667
$ws.Cells.Item(76, 11).Value = 3046.6  # K76: 3114.2856 -> 3046.6
$ws.Cells.Item(76, 12).Value = 3266.6667  # L76: 3200 -> 3266.6667
$ws.Cells.Item(76, 13).Value = -2731.6  # M76: -2799.2856 -> -2731.6
$ws.Cells.Item(76, 14).Value = -3896.6667  # N76: -3830 -> -3896.6667

$ws.Cells.Item(79, 8).Value = 3083.2778  # H79: 3145.4546 -> 3083.2778
$ws.Cells.Item(79, 9).Value = 3046.6  # I79: 3114.2856 -> 3046.6
$ws.Cells.Item(79, 10).Value = 3266.6667  # J79: 3200 -> 3266.6667
$ws.Cells.Item(79, 11).Value = 3046.6  # K79: 3114.2856 -> 3046.6
$ws.Cells.Item(79, 12).Value = 3266.6667  # L79: 3200 -> 3266.6667
$ws.Cells.Item(79, 13).Value = -1954.6  # M79: -2022.2856 -> -1954.6
$ws.Cells.Item(79, 14).Value = -5450.6667  # N79: -5384 -> -5450.6667

$ws.Cells.Item(98, 8).Value = 1234.6285  # H98: 1256.2646 -> 1234.6285
$ws.Cells.Item(98, 9).Value = 715  # I98: 725.2857 -> 715
$ws.Cells.Item(98, 11).Value = 715  # K98: 725.2857 -> 715
$ws.Cells.Item(98, 13).Value = 783  # M98: 772.7143 -> 783

$ws.Cells.Item(116, 8).Value = 9058  # H116: 4819.3125 -> 9058
$ws.Cells.Item(116, 9).Value = 1902.5  # I116: 1643.1428 -> 1902.5
$ws.Cells.Item(116, 10).Value = 11920.2  # J116: 7289.6665 -> 11920.2
$ws.Cells.Item(116, 11).Value = 1902.5  # K116: 1643.1428 -> 1902.5
$ws.Cells.Item(116, 12).Value = 11920.2  # L116: 7289.6665 -> 11920.2
$ws.Cells.Item(116, 13).Value = 1539.5  # M116: 1798.8572 -> 1539.5
$ws.Cells.Item(116, 14).Value = -18804.2  # N116: -14173.6665 -> -18804.2

$ws.Cells.Item(122, 8).Value = 1234.6285  # H122: 1256.2646 -> 1234.6285
$ws.Cells.Item(122, 9).Value = 715  # I122: 725.2857 -> 715
$ws.Cells.Item(122, 11).Value = 2145  # K122: 2175.8571 -> 2145
$ws.Cells.Item(122, 13).Value = 305  # M122: 274.1428999999998 -> 305

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5234.2295  # H32: 5408.707 -> 5234.2295
$ws.Cells.Item(32, 9).Value = 3602.04  # I32: 3680.2654 -> 3602.04
$ws.Cells.Item(32, 10).Value = 12653.272  # J32: 14819.111 -> 12653.272
$ws.Cells.Item(32, 11).Value = 3602.04  # K32: 3680.2654 -> 3602.04
$ws.Cells.Item(32, 12).Value = 12653.272  # L32: 14819.111 -> 12653.272
$ws.Cells.Item(32, 13).Value = -3315.04  # M32: -3393.2654 -> -3315.04
$ws.Cells.Item(32, 14).Value = -13227.272  # N32: -15393.111 -> -13227.272

$ws.Cells.Item(52, 8).Value = 18666.666  # H52: 18399.8 -> 18666.666
$ws.Cells.Item(52, 10).Value = 18666.666  # J52: 18399.8 -> 18666.666
$ws.Cells.Item(52, 12).Value = 18666.666  # L52: 18399.8 -> 18666.666
$ws.Cells.Item(52, 14).Value = -19302.666  # N52: -19035.8 -> -19302.666

$ws.Cells.Item(122, 8).Value = 1829.4286  # H122: 3007 -> 1829.4286
$ws.Cells.Item(122, 9).Value = 1301.3334  # I122: 0 -> 1301.3334
$ws.Cells.Item(122, 10).Value = 2780  # J122: 3007 -> 2780
$ws.Cells.Item(122, 11).Value = 3904.0002  # K122: 0 -> 3904.0002
$ws.Cells.Item(122, 12).Value = 8340  # L122: 9021 -> 8340
$ws.Cells.Item(122, 13).Value = -1454.0002  # M122: None -> -1454.0002
$ws.Cells.Item(122, 14).Value = -13240  # N122: -13921 -> -13240

$ws.Cells.Item(132, 8).Value = 2646.3584  # H132: 2894.8845 -> 2646.3584
$ws.Cells.Item(132, 9).Value = 2306.1365  # I132: 2578.325 -> 2306.1365
$ws.Cells.Item(132, 10).Value = 4309.6665  # J132: 3950.0833 -> 4309.6665
$ws.Cells.Item(132, 11).Value = 6918.4095  # K132: 7734.974999999999 -> 6918.4095
$ws.Cells.Item(132, 12).Value = 12928.9995  # L132: 11850.2499 -> 12928.9995
$ws.Cells.Item(132, 13).Value = -4388.4095  # M132: -5204.974999999999 -> -4388.4095
$ws.Cells.Item(132, 14).Value = -17988.9995  # N132: -16910.2499 -> -17988.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(42, 8).Value = 81768.336  # H42: 94000 -> 81768.336
$ws.Cells.Item(42, 9).Value = 50621  # I42: 0 -> 50621
$ws.Cells.Item(42, 10).Value = 97342  # J42: 94000 -> 97342
$ws.Cells.Item(42, 11).Value = 50621  # K42: 0 -> 50621
$ws.Cells.Item(42, 12).Value = 97342  # L42: 94000 -> 97342
$ws.Cells.Item(42, 13).Value = -50293  # M42: None -> -50293
$ws.Cells.Item(42, 14).Value = -97998  # N42: -94656 -> -97998

$ws.Cells.Item(80, 8).Value = 598.5128  # H80: 854.48 -> 598.5128
$ws.Cells.Item(80, 9).Value = 1087.8  # I80: 1282.75 -> 1087.8
$ws.Cells.Item(80, 10).Value = 292.70834  # J80: 459.15384 -> 292.70834
$ws.Cells.Item(80, 11).Value = 1087.8  # K80: 1282.75 -> 1087.8
$ws.Cells.Item(80, 12).Value = 292.70834  # L80: 459.15384 -> 292.70834
$ws.Cells.Item(80, 13).Value = -89.79999999999995  # M80: -284.75 -> -89.79999999999995
$ws.Cells.Item(80, 14).Value = -2288.70834  # N80: -2455.15384 -> -2288.70834

$ws.Cells.Item(83, 8).Value = 598.5128  # H83: 854.48 -> 598.5128
$ws.Cells.Item(83, 9).Value = 1087.8  # I83: 1282.75 -> 1087.8
$ws.Cells.Item(83, 10).Value = 292.70834  # J83: 459.15384 -> 292.70834
$ws.Cells.Item(83, 11).Value = 5439  # K83: 6413.75 -> 5439
$ws.Cells.Item(83, 12).Value = 1463.5417  # L83: 2295.7692 -> 1463.5417
$ws.Cells.Item(83, 13).Value = -447  # M83: -1421.75 -> -447
$ws.Cells.Item(83, 14).Value = -11447.5417  # N83: -12279.7692 -> -11447.5417

$ws.Cells.Item(103, 8).Value = 32750.166  # H103: 36333.668 -> 32750.166
$ws.Cells.Item(103, 10).Value = 32750.166  # J103: 36333.668 -> 32750.166
$ws.Cells.Item(103, 12).Value = 32750.166  # L103: 36333.668 -> 32750.166
$ws.Cells.Item(103, 14).Value = -35094.166  # N103: -38677.668 -> -35094.166

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 15250  # H10: 16902.334 -> 15250
$ws.Cells.Item(10, 9).Value = 500  # I10: 353.5 -> 500
$ws.Cells.Item(10, 10).Value = 30000  # J10: 50000 -> 30000
$ws.Cells.Item(10, 11).Value = 500  # K10: 353.5 -> 500
$ws.Cells.Item(10, 12).Value = 30000  # L10: 50000 -> 30000
$ws.Cells.Item(10, 13).Value = -361  # M10: -214.5 -> -361
$ws.Cells.Item(10, 14).Value = -30278  # N10: -50278 -> -30278

$ws.Cells.Item(31, 8).Value = 26317784  # H31: 27029158 -> 26317784
$ws.Cells.Item(31, 9).Value = 62500964  # I31: 71429540 -> 62500964
$ws.Cells.Item(31, 10).Value = 2742.818  # J31: 2840.9565 -> 2742.818
$ws.Cells.Item(31, 11).Value = 62500964  # K31: 71429540 -> 62500964
$ws.Cells.Item(31, 12).Value = 2742.818  # L31: 2840.9565 -> 2742.818
$ws.Cells.Item(31, 13).Value = -62500669  # M31: -71429245 -> -62500669
$ws.Cells.Item(31, 14).Value = -3332.818  # N31: -3430.9565 -> -3332.818

$ws.Cells.Item(34, 8).Value = 26317784  # H34: 27029158 -> 26317784
$ws.Cells.Item(34, 9).Value = 62500964  # I34: 71429540 -> 62500964
$ws.Cells.Item(34, 10).Value = 2742.818  # J34: 2840.9565 -> 2742.818
$ws.Cells.Item(34, 11).Value = 62500964  # K34: 71429540 -> 62500964
$ws.Cells.Item(34, 12).Value = 2742.818  # L34: 2840.9565 -> 2742.818
$ws.Cells.Item(34, 13).Value = -62500762  # M34: -71429338 -> -62500762
$ws.Cells.Item(34, 14).Value = -3146.818  # N34: -3244.9565 -> -3146.818

$ws.Cells.Item(52, 8).Value = 39000  # H52: 40000 -> 39000
$ws.Cells.Item(52, 10).Value = 39000  # J52: 40000 -> 39000
$ws.Cells.Item(52, 12).Value = 39000  # L52: 40000 -> 39000
$ws.Cells.Item(52, 14).Value = -39588  # N52: -40588 -> -39588

$ws.Cells.Item(135, 8).Value = 34243  # H135: 34212.625 -> 34243
$ws.Cells.Item(135, 10).Value = 34243  # J135: 34212.625 -> 34243
$ws.Cells.Item(135, 12).Value = 34243  # L135: 34212.625 -> 34243
$ws.Cells.Item(135, 14).Value = -44383  # N135: -44352.625 -> -44383

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 5132.727  # H56: 5079.1304 -> 5132.727
$ws.Cells.Item(56, 9).Value = 5132.727  # I56: 5079.1304 -> 5132.727
$ws.Cells.Item(56, 11).Value = 5132.727  # K56: 5079.1304 -> 5132.727
$ws.Cells.Item(56, 13).Value = -4602.727  # M56: -4549.1304 -> -4602.727

$ws.Cells.Item(124, 8).Value = 630  # H124: 4555.3335 -> 630
$ws.Cells.Item(124, 9).Value = 630  # I124: 2000 -> 630
$ws.Cells.Item(124, 10).Value = 0  # J124: 5833 -> 0
$ws.Cells.Item(124, 11).Value = 1890  # K124: 6000 -> 1890
$ws.Cells.Item(124, 12).Value = 0  # L124: 17499 -> 0
$ws.Cells.Item(124, 13).Value = 3020  # M124: -1090 -> 3020
$ws.Cells.Item(124, 14).ClearContents()  # N124: -27319 -> (removed)

$ws.Cells.Item(140, 8).Value = 3898.75  # H140: 5115.5 -> 3898.75
$ws.Cells.Item(140, 9).Value = 938  # I140: 953.3333 -> 938
$ws.Cells.Item(140, 10).Value = 8833.333000000001  # J140: 9277.666999999999 -> 8833.333000000001
$ws.Cells.Item(140, 11).Value = 2814  # K140: 2859.9999 -> 2814
$ws.Cells.Item(140, 12).Value = 26499.999  # L140: 27833.001 -> 26499.999
$ws.Cells.Item(140, 13).Value = 2366  # M140: 2320.0001 -> 2366
$ws.Cells.Item(140, 14).Value = -36859.999  # N140: -38193.001 -> -36859.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2417.1  # H132: 2339.4883 -> 2417.1
$ws.Cells.Item(132, 9).Value = 2456.6785  # I132: 2464.3215 -> 2456.6785
$ws.Cells.Item(132, 10).Value = 2324.75  # J132: 2106.4666 -> 2324.75
$ws.Cells.Item(132, 11).Value = 7370.0355  # K132: 7392.9645 -> 7370.0355
$ws.Cells.Item(132, 12).Value = 6974.25  # L132: 6319.399800000001 -> 6974.25
$ws.Cells.Item(132, 13).Value = -4840.0355  # M132: -4862.9645 -> -4840.0355
$ws.Cells.Item(132, 14).Value = -12034.25  # N132: -11379.3998 -> -12034.25

$ws.Cells.Item(137, 8).Value = 0  # H137: 57999.668 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 57999.668 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 57999.668 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137: -68199.66800000001 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(129, 8).Value = 30000  # H129: 32500 -> 30000
$ws.Cells.Item(129, 10).Value = 30000  # J129: 32500 -> 30000
$ws.Cells.Item(129, 12).Value = 30000  # L129: 32500 -> 30000
$ws.Cells.Item(129, 14).Value = -40000  # N129: -42500 -> -40000

$ws.Cells.Item(136, 8).Value = 1712  # H136: 1194.6154 -> 1712
$ws.Cells.Item(136, 9).Value = 1220.8235  # I136: 872.94116 -> 1220.8235
$ws.Cells.Item(136, 11).Value = 3662.4705  # K136: 2618.82348 -> 3662.4705
$ws.Cells.Item(136, 13).Value = -1112.4705  # M136: -68.82348000000002 -> -1112.4705

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1923.4562  # H136: 1869.339 -> 1923.4562
$ws.Cells.Item(136, 9).Value = 1702.6666  # I136: 1647.64 -> 1702.6666
$ws.Cells.Item(136, 11).Value = 5107.9998  # K136: 4942.92 -> 5107.9998
$ws.Cells.Item(136, 13).Value = -2557.9998  # M136: -2392.92 -> -2557.9998
